$wb = $excel.ActiveWorkbook

# Sheet ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1854355.4
$ws.Cells.Item(17, 10).Value = 1889330.8
$ws.Cells.Item(17, 12).Value = 5667992.4
$ws.Cells.Item(17, 14).Value = -5668328.4

# Sheet ALC row 105
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(105, 8).Value = 34892.5
$ws.Cells.Item(105, 10).Value = 34892.5
$ws.Cells.Item(105, 12).Value = 34892.5
$ws.Cells.Item(105, 14).Value = -41880.5

# Sheet ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4575.7
$ws.Cells.Item(116, 9).Value = 2122.5
$ws.Cells.Item(116, 11).Value = 2122.5
$ws.Cells.Item(116, 13).Value = 1319.5

# Sheet ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 16135060
$ws.Cells.Item(135, 9).Value = 842.05
$ws.Cells.Item(135, 10).Value = 45470000
$ws.Cells.Item(135, 11).Value = 7578.45
$ws.Cells.Item(135, 12).Value = 409230000
$ws.Cells.Item(135, 13).Value = -5043.45
$ws.Cells.Item(135, 14).Value = -409235070

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1813.2433
$ws.Cells.Item(138, 9).Value = 1152.3429
$ws.Cells.Item(138, 10).Value = 2406.359
$ws.Cells.Item(138, 11).Value = 3457.0287
$ws.Cells.Item(138, 12).Value = 7219.076999999999
$ws.Cells.Item(138, 13).Value = 1682.9713
$ws.Cells.Item(138, 14).Value = -17499.077

# Sheet ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1432.04
$ws.Cells.Item(2, 9).Value = 1512.6
$ws.Cells.Item(2, 10).Value = 1311.2
$ws.Cells.Item(2, 11).Value = 1512.6
$ws.Cells.Item(2, 12).Value = 1311.2
$ws.Cells.Item(2, 13).Value = -1399.6
$ws.Cells.Item(2, 14).Value = -1537.2

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 22852.875
$ws.Cells.Item(32, 9).Value = 22852.875
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 22852.875
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -22565.875
$ws.Cells.Item(32, 14).ClearContents()

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2882.1875
$ws.Cells.Item(45, 9).Value = 4900.2
$ws.Cells.Item(45, 10).Value = 1964.909
$ws.Cells.Item(45, 11).Value = 4900.2
$ws.Cells.Item(45, 12).Value = 1964.909
$ws.Cells.Item(45, 13).Value = -4523.2
$ws.Cells.Item(45, 14).Value = -2718.909

# Sheet ARM row 115
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(115, 8).Value = 27634.5
$ws.Cells.Item(115, 10).Value = 27634.5
$ws.Cells.Item(115, 12).Value = 27634.5
$ws.Cells.Item(115, 14).Value = -30768.5

# Sheet ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1432.04
$ws.Cells.Item(116, 9).Value = 1512.6
$ws.Cells.Item(116, 10).Value = 1311.2
$ws.Cells.Item(116, 11).Value = 1512.6
$ws.Cells.Item(116, 12).Value = 1311.2
$ws.Cells.Item(116, 13).Value = 781.4000000000001
$ws.Cells.Item(116, 14).Value = -5899.2

# Sheet BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1432.04
$ws.Cells.Item(3, 9).Value = 1512.6
$ws.Cells.Item(3, 10).Value = 1311.2
$ws.Cells.Item(3, 11).Value = 1512.6
$ws.Cells.Item(3, 12).Value = 1311.2
$ws.Cells.Item(3, 13).Value = -1398.6
$ws.Cells.Item(3, 14).Value = -1539.2

# Sheet BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1715.725
$ws.Cells.Item(86, 9).Value = 1437.1666
$ws.Cells.Item(86, 11).Value = 1437.1666
$ws.Cells.Item(86, 13).Value = -314.1666

# Sheet BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 1715.725
$ws.Cells.Item(89, 9).Value = 1437.1666
$ws.Cells.Item(89, 11).Value = 7185.833000000001
$ws.Cells.Item(89, 13).Value = -1569.833000000001

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1546.1389
$ws.Cells.Item(94, 9).Value = 1392.8928
$ws.Cells.Item(94, 10).Value = 2082.5
$ws.Cells.Item(94, 11).Value = 1392.8928
$ws.Cells.Item(94, 12).Value = 2082.5
$ws.Cells.Item(94, 13).Value = -941.8928000000001
$ws.Cells.Item(94, 14).Value = -2984.5

# Sheet BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3147.3572
$ws.Cells.Item(105, 9).Value = 3147.3572
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 3147.3572
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -1400.3572
$ws.Cells.Item(105, 14).ClearContents()

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 26895.428
$ws.Cells.Item(134, 9).Value = 31974.514
$ws.Cells.Item(134, 11).Value = 95923.542
$ws.Cells.Item(134, 13).Value = -93388.542

# Sheet BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(135, 8).Value = 48568
$ws.Cells.Item(135, 10).Value = 48568
$ws.Cells.Item(135, 12).Value = 48568
$ws.Cells.Item(135, 14).Value = -58708

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 15349.091
$ws.Cells.Item(31, 9).Value = 18875.295
$ws.Cells.Item(31, 10).Value = 3360
$ws.Cells.Item(31, 11).Value = 18875.295
$ws.Cells.Item(31, 12).Value = 3360
$ws.Cells.Item(31, 13).Value = -18580.295
$ws.Cells.Item(31, 14).Value = -3950

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 15349.091
$ws.Cells.Item(34, 9).Value = 18875.295
$ws.Cells.Item(34, 10).Value = 3360
$ws.Cells.Item(34, 11).Value = 18875.295
$ws.Cells.Item(34, 12).Value = 3360
$ws.Cells.Item(34, 13).Value = -18673.295
$ws.Cells.Item(34, 14).Value = -3764

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 14202.805
$ws.Cells.Item(132, 9).Value = 18438.69
$ws.Cells.Item(132, 11).Value = 55316.06999999999
$ws.Cells.Item(132, 13).Value = -52786.06999999999

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1114.7046
$ws.Cells.Item(134, 9).Value = 978.41174
$ws.Cells.Item(134, 11).Value = 2935.23522
$ws.Cells.Item(134, 13).Value = -400.23522

# Sheet CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 811
$ws.Cells.Item(5, 9).Value = 818.3333
$ws.Cells.Item(5, 10).Value = 800
$ws.Cells.Item(5, 11).Value = 2454.9999
$ws.Cells.Item(5, 12).Value = 2400
$ws.Cells.Item(5, 13).Value = -2342.9999
$ws.Cells.Item(5, 14).Value = -2624

# Sheet CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 597
$ws.Cells.Item(26, 9).Value = 230
$ws.Cells.Item(26, 10).Value = 719.3333
$ws.Cells.Item(26, 11).Value = 690
$ws.Cells.Item(26, 12).Value = 2157.9999
$ws.Cells.Item(26, 13).Value = -402
$ws.Cells.Item(26, 14).Value = -2733.9999

# Sheet CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 1912
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 1912
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 5736
$ws.Cells.Item(46, 14).Value = -5918
$ws.Cells.Item(46, 13).ClearContents()

# Sheet CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 610.3333
$ws.Cells.Item(92, 10).Value = 496.5
$ws.Cells.Item(92, 12).Value = 1489.5
$ws.Cells.Item(92, 14).Value = -3985.5

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 747.03
$ws.Cells.Item(131, 10).Value = 747.03
$ws.Cells.Item(131, 12).Value = 2241.09
$ws.Cells.Item(131, 14).Value = -12321.09

# Sheet CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 811
$ws.Cells.Item(135, 9).Value = 818.3333
$ws.Cells.Item(135, 10).Value = 800
$ws.Cells.Item(135, 11).Value = 7364.9997
$ws.Cells.Item(135, 12).Value = 7200
$ws.Cells.Item(135, 13).Value = -4829.9997
$ws.Cells.Item(135, 14).Value = -12270

# Sheet CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 2349.2222
$ws.Cells.Item(139, 9).Value = 1865.75
$ws.Cells.Item(139, 10).Value = 2736
$ws.Cells.Item(139, 11).Value = 5597.25
$ws.Cells.Item(139, 12).Value = 8208
$ws.Cells.Item(139, 13).Value = -457.25
$ws.Cells.Item(139, 14).Value = -18488

# Sheet CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 2470.2856
$ws.Cells.Item(140, 9).Value = 2194.9
$ws.Cells.Item(140, 11).Value = 6584.700000000001
$ws.Cells.Item(140, 13).Value = -1404.700000000001

# Sheet GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3100
$ws.Cells.Item(113, 9).Value = 2275
$ws.Cells.Item(113, 11).Value = 2275
$ws.Cells.Item(113, 13).Value = -105

# Sheet GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3045.9607
$ws.Cells.Item(126, 9).Value = 2252.5405
$ws.Cells.Item(126, 10).Value = 5142.857
$ws.Cells.Item(126, 11).Value = 6757.6215
$ws.Cells.Item(126, 12).Value = 15428.571
$ws.Cells.Item(126, 13).Value = -4287.6215
$ws.Cells.Item(126, 14).Value = -20368.571

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 59553.594
$ws.Cells.Item(132, 9).Value = 49543.137
$ws.Cells.Item(132, 10).Value = 103599.6
$ws.Cells.Item(132, 11).Value = 148629.411
$ws.Cells.Item(132, 12).Value = 310798.8
$ws.Cells.Item(132, 13).Value = -146099.411
$ws.Cells.Item(132, 14).Value = -315858.8

# Sheet LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6255.4443
$ws.Cells.Item(7, 9).Value = 6326.533
$ws.Cells.Item(7, 10).Value = 5900
$ws.Cells.Item(7, 11).Value = 6326.533
$ws.Cells.Item(7, 12).Value = 5900
$ws.Cells.Item(7, 13).Value = -6214.533
$ws.Cells.Item(7, 14).Value = -6124

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3534
$ws.Cells.Item(61, 9).Value = 2059.2
$ws.Cells.Item(61, 10).Value = 5992
$ws.Cells.Item(61, 11).Value = 2059.2
$ws.Cells.Item(61, 12).Value = 5992
$ws.Cells.Item(61, 13).Value = -1857.2
$ws.Cells.Item(61, 14).Value = -6396

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 3534
$ws.Cells.Item(113, 9).Value = 2059.2
$ws.Cells.Item(113, 10).Value = 5992
$ws.Cells.Item(113, 11).Value = 2059.2
$ws.Cells.Item(113, 12).Value = 5992
$ws.Cells.Item(113, 13).Value = 110.8000000000002
$ws.Cells.Item(113, 14).Value = -10332

# Sheet LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 6255.4443
$ws.Cells.Item(126, 9).Value = 6326.533
$ws.Cells.Item(126, 10).Value = 5900
$ws.Cells.Item(126, 11).Value = 18979.599
$ws.Cells.Item(126, 12).Value = 17700
$ws.Cells.Item(126, 13).Value = -16509.599
$ws.Cells.Item(126, 14).Value = -22640

# Sheet WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 500
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 500
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 1000
$ws.Cells.Item(100, 14).Value = -2082
$ws.Cells.Item(100, 13).ClearContents()

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1839.3334
$ws.Cells.Item(132, 9).Value = 1065.2222
$ws.Cells.Item(132, 10).Value = 3000.5
$ws.Cells.Item(132, 11).Value = 3195.6666
$ws.Cells.Item(132, 12).Value = 9001.5
$ws.Cells.Item(132, 13).Value = -665.6665999999996
$ws.Cells.Item(132, 14).Value = -14061.5

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 20834958
$ws.Cells.Item(136, 9).Value = 31251464
$ws.Cells.Item(136, 11).Value = 93754392
$ws.Cells.Item(136, 13).Value = -93751842
